$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 - match style of existing header cells (e.g. E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New boolean data columns F2:H5
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $true

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false
